$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.000",
# "0.9998", "0.000007890"). Excel's Range.Value setter auto-coerces
# such strings into numbers and silently drops the exact text
# formatting (trailing zeros, thousand-dot grouping, etc). Force
# those specific cells to Text format first so the literal string
# from the source feed is preserved exactly, matching the diff.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values (new crypto snapshot from GitHub Actions run).
$ws.Range("D2").Value = "29.971.21"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.893.37"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "0.7751"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "243.96"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.3132"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "25.81"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").Value = "0.07252"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "0.08725"
$ws.Range("E11").Value = "  +8.35%  "
$ws.Range("D12").Value = "2.051.17"
$ws.Range("E12").Value = "  +9.64%  "
$ws.Range("D13").Value = "0.7736"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "5.417"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "94.49"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").Value = "6.218"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "30.166.54"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "13.93"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "245.95"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.000007890"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.267.49"
$ws.Range("E21").Value = "  +8.60%  "
$ws.Range("D22").Value = "8.194"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "0.1608"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "9.531"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "163.19"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").Value = "18.87"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "2.052"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "1.432"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").Value = "1.545"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "4.534"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "4.133"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "1.251"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").Value = "0.7552"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "0.9998"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "2.689"
$ws.Range("E38").Value = "  +2.73%  "
$ws.Range("D39").Value = "0.01968"
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("D40").Value = "2.786"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "0.4526"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").Value = "73.63"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "1.101.08"
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("D44").Value = "6.067"
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("D45").Value = "0.8533"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "103.47"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.886"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.148.92"
$ws.Range("E49").Value = "  +6.72%  "
$ws.Range("D50").Value = "7.628"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("D51").Value = "9.884"
$ws.Range("E51").Value = "  +0.02%  "
